$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new trailing columns / header cells introduced by the edit.
$ws.Range("D1").Value = "ORG_EST_IDENOLD"
$ws.Range("E1").Value = "ORG_EST_IDENNEW"
$ws.Range("F1").Value = "ORG_EST_STATUS"

# Match the left-aligned style already used by the rest of row 1 / the sheet.
$ws.Range("D1:F1").HorizontalAlignment = -4131

# Move/collapse the selection to D9 (single cell), matching the saved view state.
$ws.Range("D9").Select() | Out-Null
